$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 159
$ws.Range("I2").Value = 150.8
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 150.8
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = -37.80000000000001
$ws.Range("N2").Value = -426

$ws.Range("H12").Value = 102.5
$ws.Range("I12").Value = 90
$ws.Range("J12").Value = 115
$ws.Range("K12").Value = 90
$ws.Range("L12").Value = 115
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = -455

$ws.Range("H132").Value = 2344
$ws.Range("I132").Value = 2503.8948
$ws.Range("J132").Value = 825
$ws.Range("K132").Value = 7511.6844
$ws.Range("L132").Value = 2475
$ws.Range("M132").Value = -4981.6844
$ws.Range("N132").Value = -7535

$ws.Range("H138").Value = 2473.85
$ws.Range("I138").Value = 1034.625
$ws.Range("J138").Value = 3151.1323
$ws.Range("K138").Value = 3103.875
$ws.Range("L138").Value = 9453.3969
$ws.Range("M138").Value = 2036.125
$ws.Range("N138").Value = -19733.3969

$ws.Range("H141").Value = 1850.4706
$ws.Range("I141").Value = 1850.4706
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5551.4118
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -371.4117999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1359.091
$ws.Range("I61").Value = 1454.5
$ws.Range("J61").Value = 1104.6666
$ws.Range("K61").Value = 1454.5
$ws.Range("L61").Value = 1104.6666
$ws.Range("M61").Value = -1242.5
$ws.Range("N61").Value = -1528.6666

$ws.Range("H74").Value = 1099.3462
$ws.Range("I74").Value = 1099.3462
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1099.3462
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -225.3462
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1099.3462
$ws.Range("I77").Value = 1099.3462
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5496.731
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1128.731
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 1084
$ws.Range("I122").Value = 968.6667
$ws.Range("J122").Value = 1257
$ws.Range("K122").Value = 2906.0001
$ws.Range("L122").Value = 3771
$ws.Range("M122").Value = -456.0001000000002
$ws.Range("N122").Value = -8671

$ws.Range("H132").Value = 2261.0588
$ws.Range("I132").Value = 1553
$ws.Range("J132").Value = 3960.4
$ws.Range("K132").Value = 4659
$ws.Range("L132").Value = 11881.2
$ws.Range("M132").Value = -2129
$ws.Range("N132").Value = -16941.2

$ws.Range("H136").Value = 1359.091
$ws.Range("I136").Value = 1454.5
$ws.Range("J136").Value = 1104.6666
$ws.Range("K136").Value = 4363.5
$ws.Range("L136").Value = 3313.9998
$ws.Range("M136").Value = -1813.5
$ws.Range("N136").Value = -8413.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 319.8
$ws.Range("I22").Value = 325.25
$ws.Range("J22").Value = 298
$ws.Range("K22").Value = 325.25
$ws.Range("L22").Value = 298
$ws.Range("M22").Value = -152.25
$ws.Range("N22").Value = -644

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1296.4546
$ws.Range("I22").Value = 1386.1
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 1386.1
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -1036.1
$ws.Range("N22").Value = -1100

$ws.Range("H122").Value = 1250
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 4031.8
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 4031.8
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 12095.4
$ws.Range("N117").Value = -18979.4
$ws.Range("M117").ClearContents()

$ws.Range("H121").Value = 724.1177
$ws.Range("I121").Value = 327.66666
$ws.Range("J121").Value = 940.36365
$ws.Range("K121").Value = 982.9999799999999
$ws.Range("L121").Value = 2821.09095
$ws.Range("M121").Value = 327.0000200000001
$ws.Range("N121").Value = -5441.09095

$ws.Range("H136").Value = 52896.95
$ws.Range("I136").Value = 112256.555
$ws.Range("J136").Value = 4330
$ws.Range("K136").Value = 336769.665
$ws.Range("L136").Value = 12990
$ws.Range("M136").Value = -331669.665
$ws.Range("N136").Value = -23190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 66.833336
$ws.Range("I2").Value = 63.75
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 63.75
$ws.Range("L2").Value = 73
$ws.Range("M2").Value = 49.25
$ws.Range("N2").Value = -299

$ws.Range("H19").Value = 15503
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 15503
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 15503
$ws.Range("N19").Value = -16079

$ws.Range("H122").Value = 775265.5600000001
$ws.Range("I122").Value = 940965.4
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2822896.2
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2820446.2
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 3423.8125
$ws.Range("I132").Value = 2917
$ws.Range("J132").Value = 4268.5
$ws.Range("K132").Value = 8751
$ws.Range("L132").Value = 12805.5
$ws.Range("M132").Value = -6221
$ws.Range("N132").Value = -17865.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 653.5
$ws.Range("I22").Value = 401
$ws.Range("J22").Value = 704
$ws.Range("K22").Value = 401
$ws.Range("L22").Value = 704
$ws.Range("M22").Value = -106
$ws.Range("N22").Value = -1294

$ws.Range("H27").Value = 653.5
$ws.Range("I27").Value = 401
$ws.Range("J27").Value = 704
$ws.Range("K27").Value = 401
$ws.Range("L27").Value = 704
$ws.Range("M27").Value = -294
$ws.Range("N27").Value = -918

$ws.Range("H33").Value = 43511.832
$ws.Range("I33").Value = 3673.3333
$ws.Range("J33").Value = 83350.336
$ws.Range("K33").Value = 3673.3333
$ws.Range("L33").Value = 83350.336
$ws.Range("M33").Value = -3383.3333
$ws.Range("N33").Value = -83930.336

$ws.Range("H122").Value = 2287.7
$ws.Range("I122").Value = 1867.7142
$ws.Range("J122").Value = 3267.6667
$ws.Range("K122").Value = 5603.142599999999
$ws.Range("L122").Value = 9803.000100000001
$ws.Range("M122").Value = -3153.142599999999
$ws.Range("N122").Value = -14703.0001

$ws.Range("H136").Value = 3402.125
$ws.Range("I136").Value = 1711.8695
$ws.Range("J136").Value = 7721.6665
$ws.Range("K136").Value = 5135.6085
$ws.Range("L136").Value = 23164.9995
$ws.Range("M136").Value = -2585.6085
$ws.Range("N136").Value = -28264.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1356.625
$ws.Range("I122").Value = 1291.3334
$ws.Range("J122").Value = 1552.5
$ws.Range("K122").Value = 3874.0002
$ws.Range("L122").Value = 4657.5
$ws.Range("M122").Value = -1424.0002
$ws.Range("N122").Value = -9557.5
